# Applies the odds corrections for Jogos_da_Semana_FlashScore_2024-11-21.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 3.7
$ws.Range("H2").Value = 3.1
$ws.Range("I2").Value = 2.15
$ws.Range("M2").Value = 1.1
$ws.Range("O2").Value = 1.5
$ws.Range("P2").Value = 2.5
$ws.Range("Q2").Value = 2.5
$ws.Range("R2").Value = 1.5

# Row 4
$ws.Range("M4").Value = 1.06
$ws.Range("N4").Value = 10
$ws.Range("O4").Value = 1.4
$ws.Range("P4").Value = 2.75
$ws.Range("Q4").Value = 2.2
$ws.Range("R4").Value = 1.65

# Row 6
$ws.Range("AI6").Value = 29
$ws.Range("AQ6").Value = 17
$ws.Range("AW6").Value = 10
$ws.Range("AZ6").Value = 251
$ws.Range("G6").Value = 1.29
$ws.Range("H6").Value = 5
$ws.Range("I6").Value = 12
$ws.Range("M6").Value = 1.04
$ws.Range("N6").Value = 13
$ws.Range("X6").Value = 6

# Row 7
$ws.Range("AD7").Value = 7.5
$ws.Range("AE7").Value = 15
$ws.Range("AK7").Value = 13
$ws.Range("AM7").Value = 201
$ws.Range("AT7").Value = 3.25
$ws.Range("AX7").Value = 9.5
$ws.Range("AZ7").Value = 29
$ws.Range("G7").Value = 4
$ws.Range("H7").Value = 3.8
$ws.Range("I7").Value = 1.8
$ws.Range("K7").Value = 2.3
$ws.Range("L7").Value = 2.4
$ws.Range("S7").Value = 1.33
$ws.Range("T7").Value = 3.25

# Row 10
$ws.Range("M10").Value = 1.1
$ws.Range("N10").Value = 7

# Row 11
$ws.Range("N11").Value = 8

# Row 14
$ws.Range("AB14").Value = 41
$ws.Range("AC14").Value = 9.5
$ws.Range("AD14").Value = 7
$ws.Range("AF14").Value = 67
$ws.Range("AH14").Value = 7.5
$ws.Range("AJ14").Value = 12
$ws.Range("AK14").Value = 15
$ws.Range("AN14").Value = 6.5
$ws.Range("AP14").Value = 34
$ws.Range("AQ14").Value = 101
$ws.Range("AR14").Value = 126
$ws.Range("AT14").Value = 2.63
$ws.Range("AW14").Value = 3.6
$ws.Range("AX14").Value = 9
$ws.Range("AZ14").Value = 29
$ws.Range("G14").Value = 4.75
$ws.Range("H14").Value = 3.5
$ws.Range("I14").Value = 1.62
$ws.Range("J14").Value = 5.5
$ws.Range("L14").Value = 2.3
$ws.Range("M14").Value = 1.06
$ws.Range("N14").Value = 10
$ws.Range("O14").Value = 1.33
$ws.Range("P14").Value = 3.25
$ws.Range("Q14").Value = 2.05
$ws.Range("R14").Value = 1.75
$ws.Range("S14").Value = 1.44
$ws.Range("T14").Value = 2.63
$ws.Range("W14").Value = 12
$ws.Range("X14").Value = 26
$ws.Range("Y14").Value = 17

# Row 15
$ws.Range("BD15").Value = 151

# Row 16
$ws.Range("M16").Value = 1.05
$ws.Range("O16").Value = 1.29
$ws.Range("Q16").Value = 1.95
$ws.Range("R16").Value = 1.9

# Row 17
$ws.Range("AH17").Value = 19
$ws.Range("AN17").Value = 4
$ws.Range("AO17").Value = 11
$ws.Range("AQ17").Value = 41
$ws.Range("AW17").Value = 5.5
$ws.Range("AZ17").Value = 67
$ws.Range("G17").Value = 1.9
$ws.Range("H17").Value = 3.25
$ws.Range("I17").Value = 3.75
$ws.Range("J17").Value = 2.63
$ws.Range("L17").Value = 4.33
$ws.Range("R17").Value = 1.75
$ws.Range("Z17").Value = 17
